$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Consignes-Input-Client")

# Insert a new row before row 3 (pushes existing rows 3.. down by one)
$ws.Rows("3:3").Insert()

# Update the text in B2 (trimestres tous régimes wording tweak)
$ws.Range("B2").Value = "trimestres :  trimestres tous régimes (pour calcul du taux de retraite et détermination carrière longue)"

# New cell C3: trimestres RG note
$ws.Range("C3").Value = "trimestres RG (trim validés pas encore utilisé dans le calcul + trim cotisés pris en compte pour calcul SAM)"

# New row 11: ATTENTION rachat note
$ws.Range("B11").Value = "ATTENTION : si rachat il faut renseigner les années concernées par le rachat (car ces années ne seront pas pris en compte dans le calcul du SAM)"

[void]$ws.Range("B12").Select()
